$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.232.49"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.499.43"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.88"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.04"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.85"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "2.891.71"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "2.489.70"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "48.109.68"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.00"
$ws.Range("E19").Value = "  +10.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.87"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "0.0₃0930"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.89"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.80"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.74"
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.65"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.37"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.59"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.84"
$ws.Range("E40").Value = "  +3.85%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.65"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "2.002.61"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.15"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.18"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.97"
$ws.Range("E51").Value = "  -1.39%  "
